$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change E3 from string "b" to the number 3
$ws.Range("E3").Value = 3

# Update the active selection to E4 (was F5)
$ws.Range("E4").Select()
